# Updated cryptos list values (price & volume columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''54.444.25'
$ws.Range('E2').Value = '  +5.31%  '
$ws.Range('D3').Value = '''3.177.71'
$ws.Range('E3').Value = '  +3.21%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''398.35'
$ws.Range('E5').Value = '  +2.66%  '
$ws.Range('D6').Value = '''109.50'
$ws.Range('E6').Value = '  +5.77%  '
$ws.Range('D7').Value = '''0.548'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.617'
$ws.Range('E9').Value = '  +4.89%  '
$ws.Range('D10').Value = '''38.98'
$ws.Range('E10').Value = '  +5.13%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = '''0.0881'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '''3.666.16'
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').Value = '''19.16'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').Value = '''8.06'
$ws.Range('E15').Value = '  +3.09%  '
$ws.Range('E16').Value = '  +8.33%  '
$ws.Range('D17').Value = '''3.180.04'
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').Value = '''10.48'
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('D19').Value = '''54.319.13'
$ws.Range('E19').Value = '  +4.88%  '
$ws.Range('D20').Value = '''3.29'
$ws.Range('E20').Value = '  +4.01%  '
$ws.Range('D21').Value = '''12.88'
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('D22').Value = '''0.0₃0988'
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('D23').Value = '''71.19'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').Value = '''272.23'
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('D26').Value = '''8.01'
$ws.Range('E26').Value = '  -2.44%  '
$ws.Range('D27').Value = '''27.70'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('D28').Value = '''7.38'
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +4.33%  '
$ws.Range('D32').Value = '''11.00'
$ws.Range('E32').Value = '  +6.79%  '
$ws.Range('D33').Value = '''0.0498'
$ws.Range('E33').Value = '  +10.64%  '
$ws.Range('D34').Value = '''36.98'
$ws.Range('E34').Value = '  +5.41%  '
$ws.Range('D35').Value = '''2.08'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').Value = '''50.55'
$ws.Range('D37').Value = '''3.62'
$ws.Range('E37').Value = '  +8.67%  '
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('E39').Value = '  +10.79%  '
$ws.Range('E40').Value = '  +9.88%  '
$ws.Range('D41').Value = '''0.292'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').Value = '''17.35'
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('D43').Value = '''1.91'
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('D44').Value = '''129.58'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').Value = '''22.21'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''2.43'
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '''2.07'
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('D49').Value = '''2.088.92'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').Value = '''0.0343'
$ws.Range('E50').Value = '  +6.94%  '
$ws.Range('B51').Value = 'FlareNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr'
$ws.Range('D51').Value = '''0.0498'
$ws.Range('E51').Value = '  +12.41%  '
